$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared text updates (header strings) ---
$ws.Range("A8").Value = "Volume 33   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/5/2026  Through  1/11/2026"

# --- Simple numeric value updates (style/format unchanged) ---
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 3
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -66.666666666666
$ws.Range("C15").Value = 7
$ws.Range("E15").Value = 250
$ws.Range("F15").Value = 20
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = 66.666666666666
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 125
$ws.Range("L15").Value = 80
$ws.Range("N15").Value = -10
$ws.Range("D16").Value = 22
$ws.Range("E16").Value = -9.090909090909
$ws.Range("F16").Value = 80
$ws.Range("G16").Value = 104
$ws.Range("H16").Value = -23.076923076923
$ws.Range("I16").Value = 32
$ws.Range("J16").Value = 37
$ws.Range("K16").Value = -13.513513513513
$ws.Range("L16").Value = 18.518518518518
$ws.Range("M16").Value = -61.44578313253
$ws.Range("N16").Value = -87.596899224806
$ws.Range("C17").Value = 51
$ws.Range("D17").Value = 44
$ws.Range("E17").Value = 15.90909090909
$ws.Range("F17").Value = 229
$ws.Range("G17").Value = 240
$ws.Range("H17").Value = -4.583333333333
$ws.Range("I17").Value = 80
$ws.Range("J17").Value = 93
$ws.Range("K17").Value = -13.978494623655
$ws.Range("L17").Value = -1.234567901234
$ws.Range("M17").Value = 77.777777777777
$ws.Range("N17").Value = -31.03448275862
$ws.Range("C18").Value = 24
$ws.Range("D18").Value = 28
$ws.Range("F18").Value = 90
$ws.Range("G18").Value = 100
$ws.Range("H18").Value = -10
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 38
$ws.Range("K18").Value = -15.78947368421
$ws.Range("L18").Value = -21.951219512195
$ws.Range("M18").Value = -45.762711864406
$ws.Range("N18").Value = -89.644012944983
$ws.Range("C19").Value = 67
$ws.Range("D19").Value = 66
$ws.Range("E19").Value = 1.515151515151
$ws.Range("F19").Value = 251
$ws.Range("G19").Value = 243
$ws.Range("H19").Value = 3.292181069958
$ws.Range("I19").Value = 93
$ws.Range("J19").Value = 92
$ws.Range("K19").Value = 1.086956521739
$ws.Range("L19").Value = 5.681818181818
$ws.Range("M19").Value = 63.157894736842
$ws.Range("N19").Value = -47.15909090909
$ws.Range("C20").Value = 40
$ws.Range("D20").Value = 42
$ws.Range("E20").Value = -4.761904761904
$ws.Range("F20").Value = 139
$ws.Range("G20").Value = 136
$ws.Range("H20").Value = 2.205882352941
$ws.Range("I20").Value = 47
$ws.Range("J20").Value = 55
$ws.Range("K20").Value = -14.545454545454
$ws.Range("L20").Value = -35.616438356164
$ws.Range("M20").Value = -14.545454545454
$ws.Range("N20").Value = -92.140468227424
$ws.Range("C21").Value = 210
$ws.Range("D21").Value = 204
$ws.Range("E21").Value = 2.941176470588
$ws.Range("F21").Value = 812
$ws.Range("G21").Value = 835
$ws.Range("H21").Value = -2.754491017964
$ws.Range("I21").Value = 294
$ws.Range("J21").Value = 319
$ws.Range("K21").Value = -7.836990595611
$ws.Range("L21").Value = -6.666666666666
$ws.Range("M21").Value = -2.649006622516
$ws.Range("N21").Value = -80
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = 12.5
$ws.Range("I22").Value = 3
$ws.Range("L22").Value = 200
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -40
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = 26.666666666666
$ws.Range("I23").Value = 6
$ws.Range("J23").Value = 6
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = -25
$ws.Range("M23").Value = 0
$ws.Range("C24").Value = 169
$ws.Range("D24").Value = 149
$ws.Range("E24").Value = 13.422818791946
$ws.Range("F24").Value = 655
$ws.Range("G24").Value = 593
$ws.Range("H24").Value = 10.455311973018
$ws.Range("I24").Value = 244
$ws.Range("J24").Value = 212
$ws.Range("K24").Value = 15.094339622641
$ws.Range("L24").Value = -2.788844621513
$ws.Range("M24").Value = 50.617283950617
$ws.Range("C25").Value = 64
$ws.Range("D25").Value = 55
$ws.Range("E25").Value = 16.363636363636
$ws.Range("F25").Value = 233
$ws.Range("G25").Value = 207
$ws.Range("H25").Value = 12.56038647343
$ws.Range("I25").Value = 80
$ws.Range("J25").Value = 79
$ws.Range("K25").Value = 1.26582278481
$ws.Range("L25").Value = -20.79207920792
$ws.Range("C26").Value = 67
$ws.Range("D26").Value = 71
$ws.Range("E26").Value = -5.633802816901
$ws.Range("F26").Value = 333
$ws.Range("G26").Value = 344
$ws.Range("H26").Value = -3.197674418604
$ws.Range("I26").Value = 120
$ws.Range("J26").Value = 126
$ws.Range("K26").Value = -4.761904761904
$ws.Range("L26").Value = -7.692307692307
$ws.Range("M26").Value = -6.25
$ws.Range("C27").Value = 8
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 21
$ws.Range("G27").Value = 14
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 10
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = 66.666666666666
$ws.Range("L27").Value = 42.857142857142
$ws.Range("C28").Value = 11
$ws.Range("D28").Value = 7
$ws.Range("E28").Value = 57.142857142857
$ws.Range("F28").Value = 38
$ws.Range("G28").Value = 35
$ws.Range("H28").Value = 8.571428571428
$ws.Range("I28").Value = 18
$ws.Range("J28").Value = 15
$ws.Range("K28").Value = 20
$ws.Range("L28").Value = 38.461538461538
$ws.Range("G31").Value = 4
$ws.Range("H31").Value = -75
$ws.Range("J31").Value = 2
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 3
$ws.Range("H33").Value = -66.666666666666
$ws.Range("J43").Value = 1080
$ws.Range("K43").Value = -69.508752117447
$ws.Range("L43").Value = -78.984238178634
$ws.Range("M43").Value = -88.40455228688
$ws.Range("N43").Value = -90.981963927855
$ws.Range("J45").Value = 1933
$ws.Range("K45").Value = -57.206110250166
$ws.Range("L45").Value = -73.374655647382
$ws.Range("M45").Value = -90.021165660007
$ws.Range("N45").Value = -91.115911388914
$ws.Range("J46").Value = 11122
$ws.Range("K46").Value = -33.872406207265
$ws.Range("L46").Value = -53.474168583978
$ws.Range("M46").Value = -77.400280413711
$ws.Range("N46").Value = -79.43759359574

# --- Cells changing FROM placeholder text TO a number (need numeric style) ---
# I14, J22: text "0" -> number, style 13 -> 14 (#,##0)
$ws.Range("C14").Copy()
$ws.Range("I14,J22").PasteSpecial(-4122)
$ws.Range("I14").Value = 1
$ws.Range("J22").Value = 1

# M15, K22, M22, L29, L30: text "***.*" -> number, style 13 -> 15 (#,##0.0)
$ws.Range("L22").Copy()
$ws.Range("M15,K22,M22,L29,L30").PasteSpecial(-4122)
$ws.Range("M15").Value = 350
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = 200
$ws.Range("L29").Value = -100
$ws.Range("L30").Value = -100

# --- Cells changing FROM a number TO placeholder text (need General/text style) ---
# C29, D29, C30, D30, C33: number -> text "0", style 14 -> 13 (General)
$ws.Range("C29").Value = "'0"
$ws.Range("D29").Value = "'0"
$ws.Range("C30").Value = "'0"
$ws.Range("D30").Value = "'0"
$ws.Range("C33").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("C29,D29,C30,D30,C33").PasteSpecial(-4122)

# E29, E30: number -> text "***.*", style 15 -> 13 (General)
$ws.Range("E29").Value = "'***.*"
$ws.Range("E30").Value = "'***.*"
$ws.Range("D14").Copy()
$ws.Range("E29,E30").PasteSpecial(-4122)

